# Testprotokoll.xlsx edit:
# - Insert a new row (row 5) with a new task item
#   "Transaktionen korrekt verschlüsseln"
# - Highlight column B next to the first 4 rows (header + the 3 "austesten"
#   rows) with a solid green fill
# - Leave the final selection on K15 (matches the author's last click)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row above the current row 5 ("Transaktionen validieren"),
# pushing everything below it down by one row, then fill in the new cell.
$ws.Rows(5).Insert()
$ws.Range("A5").Value = "Transaktionen korrekt verschlüsseln"

# Apply a solid green fill (RGB 0,176,80 -> FF00B050) to B1:B4.
$ws.Range("B1:B4").Interior.Color = 5287936

# Restore the selection to the cell that was active when the workbook was
# last saved.
$null = $ws.Range("K15").Select()
